# Remove the duplicate, visible "R.C. 2943.031." reference (and its
# leading space) that immediately follows the citizenship-advisement
# sentence, while leaving the already-hidden "R.C. 2943.031." runs
# further along untouched.
#
# Before:  ...under United States law.[ ][R.C. 2943.031.][ ][hidden: R.C. 2943.031.][ ]The Court, finding...
# After:   ...under United States law.[ ][hidden: R.C. 2943.031.][ ]The Court, finding...

$d = $word.ActiveDocument

# Locate the end of the sentence that precedes the duplicate citation.
$anchorEnd = $d.Content
$foundEnd = $anchorEnd.Find.Execute(
    "under United States law.", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $foundEnd) {
    throw "Could not locate the 'under United States law.' anchor text."
}
$lawEnd = $anchorEnd.End

# Locate the start of the following sentence ("The Court, finding...").
$anchorStart = $d.Content
$foundStart = $anchorStart.Find.Execute(
    "The Court, finding that the Defendant entered the plea knowingly",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundStart) {
    throw "Could not locate the 'The Court, finding...' anchor text."
}
$courtStart = $anchorStart.Start

# Everything visible between the two anchors is " R.C. 2943.031. " (a
# leading space, the duplicate citation, then a trailing space). Keep
# the trailing space (it precedes hidden text that must survive), and
# delete only the leading space + duplicate citation.
$dupRange = $d.Range($lawEnd, $courtStart - 1)

if ($dupRange.Text -ne " R.C. 2943.031.") {
    throw "Unexpected text in range to delete: [$($dupRange.Text)]"
}

$dupRange.Delete()
